$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 208.4
$ws.Range("I6").Value = 208.4
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 625.2
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -513.2
$ws.Range("N6").ClearContents()

$ws.Range("H12").Value = 4996
$ws.Range("I12").Value = 4996
$ws.Range("K12").Value = 4996
$ws.Range("M12").Value = -4826

$ws.Range("H21").Value = 3165.6667
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

$ws.Range("H23").Value = 3165.6667
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H29").Value = 290.4
$ws.Range("I29").Value = 290.4
$ws.Range("K29").Value = 871.1999999999999
$ws.Range("M29").Value = -590.1999999999999

$ws.Range("H32").Value = 7100
$ws.Range("J32").Value = 7250
$ws.Range("L32").Value = 7250
$ws.Range("N32").Value = -7902

$ws.Range("H43").Value = 3999.75
$ws.Range("I43").Value = 2500
$ws.Range("K43").Value = 2500
$ws.Range("M43").Value = -2431

$ws.Range("H58").Value = 837.5
$ws.Range("J58").Value = 1000
$ws.Range("L58").Value = 3000
$ws.Range("N58").Value = -3300

$ws.Range("H92").Value = 503.05884
$ws.Range("J92").Value = 233.28572
$ws.Range("L92").Value = 233.28572
$ws.Range("N92").Value = -2729.28572

$ws.Range("H98").Value = 832.2778
$ws.Range("I98").Value = 735.1429000000001
$ws.Range("J98").Value = 1172.25
$ws.Range("K98").Value = 735.1429000000001
$ws.Range("L98").Value = 1172.25
$ws.Range("M98").Value = 762.8570999999999
$ws.Range("N98").Value = -4168.25

$ws.Range("H107").Value = 1755
$ws.Range("I107").Value = 1755
$ws.Range("K107").Value = 1755
$ws.Range("M107").Value = 165

$ws.Range("H122").Value = 832.2778
$ws.Range("I122").Value = 735.1429000000001
$ws.Range("J122").Value = 1172.25
$ws.Range("K122").Value = 2205.4287
$ws.Range("L122").Value = 3516.75
$ws.Range("M122").Value = 244.5712999999996
$ws.Range("N122").Value = -8416.75

$ws.Range("H131").Value = 2396.3333
$ws.Range("I131").Value = 2396.3333
$ws.Range("K131").Value = 7188.999899999999
$ws.Range("M131").Value = -2148.999899999999

$ws.Range("H137").Value = 1378
$ws.Range("I137").Value = 1327.7142
$ws.Range("J137").Value = 1428.2858
$ws.Range("K137").Value = 3983.1426
$ws.Range("L137").Value = 4284.857400000001
$ws.Range("M137").Value = -1433.1426
$ws.Range("N137").Value = -9384.857400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()

$ws.Range("H37").Value = 8123.75
$ws.Range("I37").Value = 5000
$ws.Range("J37").Value = 9998
$ws.Range("K37").Value = 5000
$ws.Range("L37").Value = 9998
$ws.Range("M37").Value = -4727
$ws.Range("N37").Value = -10544

$ws.Range("H61").Value = 2032.2858
$ws.Range("J61").Value = 1778.6666
$ws.Range("L61").Value = 1778.6666
$ws.Range("N61").Value = -2202.6666

$ws.Range("H110").Value = 250
$ws.Range("I110").Value = 250
$ws.Range("K110").Value = 250
$ws.Range("M110").Value = 1795

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws.Range("H136").Value = 2032.2858
$ws.Range("J136").Value = 1778.6666
$ws.Range("L136").Value = 5335.9998
$ws.Range("N136").Value = -10435.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 195.42857
$ws.Range("I107").Value = 195
$ws.Range("K107").Value = 195
$ws.Range("M107").Value = 1725

$ws.Range("H134").Value = 6851.1
$ws.Range("I134").Value = 6556.778
$ws.Range("J134").Value = 9500
$ws.Range("K134").Value = 19670.334
$ws.Range("L134").Value = 28500
$ws.Range("M134").Value = -17135.334
$ws.Range("N134").Value = -33570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 483.75
$ws.Range("I22").Value = 467.5
$ws.Range("K22").Value = 467.5
$ws.Range("M22").Value = -117.5

$ws.Range("H33").Value = 10004
$ws.Range("I33").Value = 10004
$ws.Range("K33").Value = 10004
$ws.Range("M33").Value = -9625

$ws.Range("H107").Value = 821.1429000000001
$ws.Range("I107").Value = 649.6
$ws.Range("J107").Value = 1250
$ws.Range("K107").Value = 649.6
$ws.Range("L107").Value = 1250
$ws.Range("M107").Value = 1270.4
$ws.Range("N107").Value = -5090

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1235
$ws.Range("I5").Value = 1266.75
$ws.Range("J5").Value = 1216.8572
$ws.Range("K5").Value = 3800.25
$ws.Range("L5").Value = 3650.5716
$ws.Range("M5").Value = -3688.25
$ws.Range("N5").Value = -3874.5716

$ws.Range("H17").Value = 699
$ws.Range("J17").Value = 699
$ws.Range("L17").Value = 2097
$ws.Range("N17").Value = -2435

$ws.Range("H34").Value = 212.5
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws.Range("H39").Value = 640
$ws.Range("J39").Value = 1200
$ws.Range("L39").Value = 3600
$ws.Range("N39").Value = -4188

$ws.Range("H55").Value = 1037.6
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws.Range("H135").Value = 1235
$ws.Range("I135").Value = 1266.75
$ws.Range("J135").Value = 1216.8572
$ws.Range("K135").Value = 11400.75
$ws.Range("L135").Value = 10951.7148
$ws.Range("M135").Value = -8865.75
$ws.Range("N135").Value = -16021.7148

$ws.Range("H137").Value = 3000
$ws.Range("I137").Value = 3000
$ws.Range("K137").Value = 9000
$ws.Range("M137").Value = -3900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 268.27777
$ws.Range("I2").Value = 57.77778
$ws.Range("J2").Value = 478.77777
$ws.Range("K2").Value = 57.77778
$ws.Range("L2").Value = 478.77777
$ws.Range("M2").Value = 55.22222
$ws.Range("N2").Value = -704.7777699999999

$ws.Range("H122").Value = 333929.66
$ws.Range("I122").Value = 333929.66
$ws.Range("K122").Value = 1001788.98
$ws.Range("M122").Value = -999338.98

$ws.Range("H132").Value = 4194
$ws.Range("I132").Value = 4014.3635
$ws.Range("J132").Value = 4758.5713
$ws.Range("K132").Value = 12043.0905
$ws.Range("L132").Value = 14275.7139
$ws.Range("M132").Value = -9513.0905
$ws.Range("N132").Value = -19335.7139

$ws.Range("H139").Value = 90000
$ws.Range("J139").Value = 90000
$ws.Range("L139").Value = 90000
$ws.Range("N139").Value = -100280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 6959.6
$ws.Range("I16").Value = 6959.6
$ws.Range("K16").Value = 6959.6
$ws.Range("M16").Value = -6789.6

$ws.Range("H46").Value = 3673.2666
$ws.Range("I46").Value = 3511
$ws.Range("J46").Value = 3916.6667
$ws.Range("K46").Value = 3511
$ws.Range("L46").Value = 3916.6667
$ws.Range("M46").Value = -3323
$ws.Range("N46").Value = -4292.6667

$ws.Range("I55").Value = 975.25
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 975.25
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -802.25
$ws.Range("N55").ClearContents()

$ws.Range("H93").Value = 918
$ws.Range("I93").Value = 918
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 918
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 330
$ws.Range("N93").ClearContents()

$ws.Range("H101").Value = 13575.667
$ws.Range("J101").Value = 13575.667
$ws.Range("L101").Value = 13575.667
$ws.Range("N101").Value = -20065.667

$ws.Range("H122").Value = 4945
$ws.Range("I122").Value = 4564.5
$ws.Range("K122").Value = 13693.5
$ws.Range("M122").Value = -11243.5

$ws.Range("H132").Value = 7003
$ws.Range("J132").Value = 8024.2
$ws.Range("L132").Value = 24072.6
$ws.Range("N132").Value = -29132.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1619.75
$ws.Range("J107").Value = 1600
$ws.Range("L107").Value = 4800
$ws.Range("N107").Value = -8640

$ws.Range("H132").Value = 1552.5
$ws.Range("I132").Value = 1803.5555
$ws.Range("J132").Value = 799.3333
$ws.Range("K132").Value = 5410.666499999999
$ws.Range("L132").Value = 2397.9999
$ws.Range("M132").Value = -2880.666499999999
$ws.Range("N132").Value = -7457.9999
